$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Shorten the "For the game to work..." paragraph: drop the tail
#    about $random/$urandom/counter/16-bit register (that content
#    moves into a new paragraph below, reworded).
# ------------------------------------------------------------------
$d.Content.Find.Execute(
    ", but instead of using the `$random or `$urandom commands, I instead used a counter that counts at the speed of the 100Mhz clock. This number is a 16 bit register, so we can use it as 4 random 4 bit numbers. ",
    $true, $false, $false, $false, $false, $true, 1, $false, ". ", 2) | Out-Null

# ------------------------------------------------------------------
# 2. Re-sequence the "B. Implementation" heading so it appears right
#    before the (still-to-come) "We use this 16 bit number..." text.
#    Simplest reliable approach: swap the *text* of the two existing
#    paragraphs that currently hold "We use this 16 bit number..." and
#    "B. Implementation", then insert a brand-new paragraph between
#    them carrying the new "Instead of using..." text.
# ------------------------------------------------------------------
$weUseText = "We use this 16 bit number to basically choose which output will be on the map. We decided that we wanted 3 different types of outputs for each section of the map: empty, bottom block, and top block. This means we need 2 bits to represent each spot on the map. There are 8 segmented displays on the fpga, so we need 16 bits to represent the map. "

$pWeUse = $d.Paragraphs.Item(4)
$pImpl  = $d.Paragraphs.Item(5)

$pWeUse.Range.Text = "B. Implementation"

$pImpl.Range.InsertParagraphAfter()
$pNew = $d.Paragraphs.Item(6)
$pNew.Range.Text = $weUseText

$pImpl.Range.Text = "Instead of using the `$random or `$urandom commands, I instead used a counter that counts at the speed of the 100Mhz clock. This number is a 16 bit register, so we can use it as 4 random 4 bit numbers. Note that, although the register is named seed, it is not actually acting as a seed, it’s just a random number."
$pImpl.LeftIndent = 0

# ------------------------------------------------------------------
# 3. Update the "To break down the Always block..." paragraph text.
# ------------------------------------------------------------------
$d.Content.Find.Execute(
    "To break down the Always block that starts on line 94, every time",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "To break down the Always block that starts on line 94 (the largest part of the program), every time", 2) | Out-Null

$d.Content.Find.Execute(
    "and when there is not any preload left, the portion of the seed the program is currently on is what determines what the next preload is. ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "and when there is not any preload left, the seed determines what the next preload is. ", 2) | Out-Null

# ------------------------------------------------------------------
# 4. Clear the trailing tab-only paragraph and append a new
#    bibliography paragraph after it.
# ------------------------------------------------------------------
$pTab = $d.Paragraphs.Last
$pTab.Range.Text = ""

$pTab.Range.InsertParagraphAfter()
$pBib = $d.Paragraphs.Last
$pBib.Range.Text = "For bibliography: https://www.eecs.umich.edu/courses/eecs270/270lab/270_docs/debounce.html"

Write-Output "done"
